$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Sheet is protected; unprotect before editing, re-protect afterwards
$ws.Unprotect()

# Update the confidential/model-holdings date string in cell A38
$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.
Model holdings provided as of 2021-05-06 for illustrative purposes only and are subject to change."

# Update the numeric holdings data (columns D and E) for rows 2-35
$updates = @(
    @{ Row = 2; D = 0.03620584348821039; E = 0 },
    @{ Row = 3; D = 0.02046683009110789; E = 0.001168679392286487 },
    @{ Row = 4; D = 0.01923664916165584; E = 0.001187303125757211 },
    @{ Row = 5; D = 0.03783214210452159; E = 0.0003512469265893881 },
    @{ Row = 6; D = 0.03434352202431512; E = 0.0004003202562048447 },
    @{ Row = 7; D = 0.0198751374738535; E = -0.0006750241080039077 },
    @{ Row = 8; D = 0.03716184335050143; E = 0.004903005755702461 },
    @{ Row = 9; D = 0.02040304565763346; E = 0.001893768599513068 },
    @{ Row = 10; D = 0.02539111101770599; E = 0.01187600644122377 },
    @{ Row = 11; D = 0.02381055524431115; E = 0.01046418030587604 },
    @{ Row = 12; D = 0.05724576240347597; E = 0.009535160905840279 },
    @{ Row = 13; D = 0.02513352003636694; E = -0.001830161054172774 },
    @{ Row = 14; D = 0.02742986186008983; E = 0.00916730328495019 },
    @{ Row = 15; D = 0.03359558821064139; E = 0.01022320667916166 },
    @{ Row = 16; D = 0.01981554400475799; E = -0.01951463597698266 },
    @{ Row = 17; D = 0.03015961073178; E = 0.02491103202846978 },
    @{ Row = 18; D = 0.04218358975357048; E = -0.0009208103130755596 },
    @{ Row = 19; D = 0.1265536313827263; E = -0.001326259946949571 },
    @{ Row = 20; D = 0.009133746879980349; E = 0.001029601029600924 },
    @{ Row = 21; D = 0.01533739634722905; E = 0.008467459762071528 },
    @{ Row = 22; D = 0.01633106377960876; E = 0.0283570972547349 },
    @{ Row = 23; D = 0.01564313230960408; E = 0.004273504273504258 },
    @{ Row = 24; D = 0.02125396475452111; E = 0.005376900010340302 },
    @{ Row = 25; D = 0.01239605488372466; E = 0.006629834254143541 },
    @{ Row = 26; D = 0.04184473495073757; E = 0.01330108827085863 },
    @{ Row = 27; D = 0.02396638756615693; E = 0 },
    @{ Row = 28; D = 0.04565708147455545; E = 0.001430615164520921 },
    @{ Row = 29; D = 0.05548652844143801; E = 0.01245712222422823 },
    @{ Row = 30; D = 0.01322678579854308; E = 0.0006414368184732844 },
    @{ Row = 31; D = 0.02067627609141095; E = 0.001915708812260553 },
    @{ Row = 32; D = 0.01347369493164803; E = 0.005215742057847361 },
    @{ Row = 33; D = 0.0420278085410464; E = -0.0005144032921812203 },
    @{ Row = 34; D = 0.01670155525257041; E = 0.007497375918428473 },
    @{ Row = 35; D = $null; E = 0.004485373493430256 }

)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Cells.Item($u.Row, 4).Value = $u.D
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# Re-apply sheet protection (sheet was protected before this edit)
$ws.Protect()
